$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells P1:Q1 - copy style (bold + border) from O1, then set values
$ws.Range("O1").Copy($ws.Range("P1:Q1"))
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

$row2 = New-Object "object[,]" 1,16
$row2[0,0] = 24.03951969628692
$row2[0,1] = 0
$row2[0,2] = 2.072646856124152
$row2[0,3] = 31.51251701979448
$row2[0,4] = 24.56165965671888
$row2[0,5] = 35.9267973137101
$row2[0,6] = 4.413919534380266
$row2[0,7] = 3.005390845382122
$row2[0,8] = 9.83149814950481
$row2[0,9] = 70.04831980844861
$row2[0,10] = 0
$row2[0,11] = 0
$row2[0,12] = 0
$row2[0,13] = 0
$row2[0,14] = 0
$row2[0,15] = 18.80947629602494
$ws.Range("B2:Q2").Value = $row2

$row3 = New-Object "object[,]" 1,16
$row3[0,0] = 22.47103589359181
$row3[0,1] = 0
$row3[0,2] = 2.011273384765137
$row3[0,3] = 29.62929689500078
$row3[0,4] = 23.04066960629008
$row3[0,5] = 33.49013665598645
$row3[0,6] = 4.141136017262474
$row3[0,7] = 3.192510874708861
$row3[0,8] = 9.399448728397596
$row3[0,9] = 65.44940052185392
$row3[0,10] = 0
$row3[0,11] = 0
$row3[0,12] = 0
$row3[0,13] = 0
$row3[0,14] = 0
$row3[0,15] = 17.7249640292932
$ws.Range("B3:Q3").Value = $row3

$row4 = New-Object "object[,]" 1,16
$row4[0,0] = 21.44996014040526
$row4[0,1] = 0
$row4[0,2] = 1.971595592722466
$row4[0,3] = 28.41373953811855
$row4[0,4] = 22.0654961429624
$row4[0,5] = 31.93724424972372
$row4[0,6] = 3.968141076107703
$row4[0,7] = 3.308238040789099
$row4[0,8] = 9.128193636010886
$row4[0,9] = 62.55294165479136
$row4[0,10] = 0
$row4[0,11] = 0
$row4[0,12] = 0
$row4[0,13] = 0
$row4[0,14] = 0
$row4[0,15] = 17.0328935374161
$ws.Range("B4:Q4").Value = $row4

$row5 = New-Object "object[,]" 1,16
$row5[0,0] = 21.0191248006034
$row5[0,1] = 0
$row5[0,2] = 1.954931188485111
$row5[0,3] = 27.90356262375301
$row5[0,4] = 21.65793092303546
$row5[0,5] = 31.28541122636794
$row5[0,6] = 3.89622832676626
$row5[0,7] = 3.357057364133082
$row5[0,8] = 9.016283131538131
$row5[0,9] = 61.41579519922216
$row5[0,10] = 0
$row5[0,11] = 0
$row5[0,12] = 0
$row5[0,13] = 0
$row5[0,14] = 0
$row5[0,15] = 16.7444962881794
$ws.Range("B5:Q5").Value = $row5

$row6 = New-Object "object[,]" 1,16
$row6[0,0] = 20.94660719316666
$row6[0,1] = 0
$row6[0,2] = 1.952134346290389
$row6[0,3] = 27.81795549951648
$row6[0,4] = 21.58964784109698
$row6[0,5] = 31.17602849022096
$row6[0,6] = 3.884200850900155
$row6[0,7] = 3.367086934272987
$row6[0,8] = 8.997621809765118
$row6[0,9] = 61.31282198545149
$row6[0,10] = 0
$row6[0,11] = 0
$row6[0,12] = 0
$row6[0,13] = 0
$row6[0,14] = 0
$row6[0,15] = 16.6962303968986
$ws.Range("B6:Q6").Value = $row6

$row7 = New-Object "object[,]" 1,16
$row7[0,0] = 21.44400016915299
$row7[0,1] = 0
$row7[0,2] = 1.971372848072053
$row7[0,3] = 28.40691902513552
$row7[0,4] = 22.06004045089546
$row7[0,5] = 31.92853030908462
$row7[0,6] = 3.96717700955433
$row7[0,7] = 3.31378638727711
$row7[0,8] = 9.126689739862357
$row7[0,9] = 62.77205030764391
$row7[0,10] = 0
$row7[0,11] = 0
$row7[0,12] = 0
$row7[0,13] = 0
$row7[0,14] = 0
$row7[0,15] = 17.02902959968149
$ws.Range("B7:Q7").Value = $row7

$row8 = New-Object "object[,]" 1,16
$row8[0,0] = 23.51060482620635
$row8[0,1] = 0
$row8[0,2] = 2.051902659387873
$row8[0,3] = 30.87551171780761
$row8[0,4] = 24.04589402848931
$row8[0,5] = 35.09754434051205
$row8[0,6] = 4.320986758575229
$row8[0,7] = 3.075717945874242
$row8[0,8] = 9.683840104137671
$row8[0,9] = 68.7835940212962
$row8[0,10] = 0
$row8[0,11] = 0
$row8[0,12] = 0
$row8[0,13] = 0
$row8[0,14] = 0
$row8[0,15] = 18.44107878545448
$ws.Range("B8:Q8").Value = $row8

$row9 = New-Object "object[,]" 1,16
$row9[0,0] = 27.11809087931507
$row9[0,1] = 0
$row9[0,2] = 2.193758469364389
$row9[0,3] = 35.2450495340266
$row9[0,4] = 27.62587036344699
$row9[0,5] = 40.77835512480394
$row9[0,6] = 4.972639787590964
$row9[0,7] = 2.607337829569494
$row9[0,8] = 10.72423450015806
$row9[0,9] = 79.22058890756904
$row9[0,10] = 0
$row9[0,11] = 0
$row9[0,12] = 0
$row9[0,13] = 0
$row9[0,14] = 0
$row9[0,15] = 20.99573808943829
$ws.Range("B9:Q9").Value = $row9

$row10 = New-Object "object[,]" 1,16
$row10[0,0] = 29.52590840198576
$row10[0,1] = 0
$row10[0,2] = 2.310072501018853
$row10[0,3] = 37.29635707322597
$row10[0,4] = 30.03866108906345
$row10[0,5] = 44.30907731206317
$row10[0,6] = 5.385526929647921
$row10[0,7] = 2.694300863279516
$row10[0,8] = 11.39247492535364
$row10[0,9] = 86.33114078562375
$row10[0,10] = 0
$row10[0,11] = 0
$row10[0,12] = 0
$row10[0,13] = 0
$row10[0,14] = 0
$row10[0,15] = 22.75496679192758
$ws.Range("B10:Q10").Value = $row10

$row11 = New-Object "object[,]" 1,16
$row11[0,0] = 30.56501963939844
$row11[0,1] = 0
$row11[0,2] = 2.563663453215443
$row11[0,3] = 30.81223697021244
$row11[0,4] = 29.78429903360374
$row11[0,5] = 43.58311598726675
$row11[0,6] = 5.616659145951649
$row11[0,7] = 2.863409948912086
$row11[0,8] = 11.17674373362729
$row11[0,9] = 89.71136125585014
$row11[0,10] = 0
$row11[0,11] = 0
$row11[0,12] = 0
$row11[0,13] = 0
$row11[0,14] = 0
$row11[0,15] = 22.39713327810053
$ws.Range("B11:Q11").Value = $row11

$row12 = New-Object "object[,]" 1,16
$row12[0,0] = 30.94890444533592
$row12[0,1] = 0
$row12[0,2] = 2.767133435402345
$row12[0,3] = 25.05731502740487
$row12[0,4] = 29.0754691611019
$row12[0,5] = 42.24736429732412
$row12[0,6] = 6.237937671084402
$row12[0,7] = 2.920602483785713
$row12[0,8] = 10.84850170145483
$row12[0,9] = 90.8143682008233
$row12[0,10] = 0
$row12[0,11] = 0
$row12[0,12] = 0
$row12[0,13] = 0
$row12[0,14] = 0
$row12[0,15] = 21.72326754345866
$ws.Range("B12:Q12").Value = $row12

$row13 = New-Object "object[,]" 1,16
$row13[0,0] = 30.86632728123295
$row13[0,1] = 0
$row13[0,2] = 2.943809514685327
$row13[0,3] = 19.49624690559535
$row13[0,4] = 27.95062371099606
$row13[0,5] = 40.31767639658085
$row13[0,6] = 7.085888880659661
$row13[0,7] = 2.906780731808359
$row13[0,8] = 10.40524599824509
$row13[0,9] = 90.546897089683
$row13[0,10] = 0
$row13[0,11] = 0
$row13[0,12] = 0
$row13[0,13] = 0
$row13[0,14] = 0
$row13[0,15] = 20.74356323316852
$ws.Range("B13:Q13").Value = $row13

$row14 = New-Object "object[,]" 1,16
$row14[0,0] = 30.59665531019523
$row14[0,1] = 0
$row14[0,2] = 3.056828410694433
$row14[0,3] = 15.85806279836412
$row14[0,4] = 26.95979231451662
$row14[0,5] = 38.67195191279288
$row14[0,6] = 7.782501013475196
$row14[0,7] = 2.867465423457775
$row14[0,8] = 10.03791601949647
$row14[0,9] = 89.78932321383597
$row14[0,10] = 0
$row14[0,11] = 0
$row14[0,12] = 0
$row14[0,13] = 0
$row14[0,14] = 0
$row14[0,15] = 19.90649323966064
$ws.Range("B14:Q14").Value = $row14

$row15 = New-Object "object[,]" 1,16
$row15[0,0] = 30.4309102586133
$row15[0,1] = 0
$row15[0,2] = 3.077355927417608
$row15[0,3] = 15.01738713704774
$row15[0,4] = 26.62442171443513
$row15[0,5] = 38.1354380645318
$row15[0,6] = 7.943310175472848
$row15[0,7] = 2.846492030909852
$row15[0,8] = 9.923695739391185
$row15[0,9] = 89.384250240553
$row15[0,10] = 0
$row15[0,11] = 0
$row15[0,12] = 0
$row15[0,13] = 0
$row15[0,14] = 0
$row15[0,15] = 19.63251974898982
$ws.Range("B15:Q15").Value = $row15

$row16 = New-Object "object[,]" 1,16
$row16[0,0] = 29.46102307622637
$row16[0,1] = 0
$row16[0,2] = 2.997938480660388
$row16[0,3] = 14.87812078050915
$row16[0,4] = 25.76598784452766
$row16[0,5] = 36.91053950359841
$row16[0,6] = 7.660739409246655
$row16[0,7] = 2.708061172222736
$row16[0,8] = 9.708386443791598
$row16[0,9] = 86.6367325381926
$row16[0,10] = 0
$row16[0,11] = 0
$row16[0,12] = 0
$row16[0,13] = 0
$row16[0,14] = 0
$row16[0,15] = 18.99672021100486
$ws.Range("B16:Q16").Value = $row16

$row17 = New-Object "object[,]" 1,16
$row17[0,0] = 28.85101104749148
$row17[0,1] = 0
$row17[0,2] = 2.873301182057784
$row17[0,3] = 16.67373059312068
$row17[0,4] = 25.65108559226206
$row17[0,5] = 36.87753641026156
$row17[0,6] = 7.024769438298595
$row17[0,7] = 2.623905086322607
$row17[0,8] = 9.745272924165569
$row17[0,9] = 84.91638599100503
$row17[0,10] = 0
$row17[0,11] = 0
$row17[0,12] = 0
$row17[0,13] = 0
$row17[0,14] = 0
$row17[0,15] = 18.969505310182
$ws.Range("B17:Q17").Value = $row17

$row18 = New-Object "object[,]" 1,16
$row18[0,0] = 28.49337206267155
$row18[0,1] = 0
$row18[0,2] = 2.700110207342719
$row18[0,3] = 20.80424627189185
$row18[0,4] = 26.17968766275382
$row18[0,5] = 37.89277141949445
$row18[0,6] = 6.12443301700185
$row18[0,7] = 2.567073126085487
$row18[0,8] = 10.00238808494644
$row18[0,9] = 83.74900462306361
$row18[0,10] = 0
$row18[0,11] = 0
$row18[0,12] = 0
$row18[0,13] = 0
$row18[0,14] = 0
$row18[0,15] = 19.4779851401496
$ws.Range("B18:Q18").Value = $row18

$row19 = New-Object "object[,]" 1,16
$row19[0,0] = 28.37226976796862
$row19[0,1] = 0
$row19[0,2] = 2.508970537997239
$row19[0,3] = 26.76589938003147
$row19[0,4] = 27.18245948931966
$row19[0,5] = 39.66748395264133
$row19[0,6] = 5.355079688917771
$row19[0,7] = 2.554485117098951
$row19[0,8] = 10.41147892583367
$row19[0,9] = 83.47852690970711
$row19[0,10] = 0
$row19[0,11] = 0
$row19[0,12] = 0
$row19[0,13] = 0
$row19[0,14] = 0
$row19[0,15] = 20.3773671542423
$ws.Range("B19:Q19").Value = $row19

$row20 = New-Object "object[,]" 1,16
$row20[0,0] = 28.91642397048742
$row20[0,1] = 0
$row20[0,2] = 2.281179092170296
$row20[0,3] = 36.72981096855779
$row20[0,4] = 29.41374269254153
$row20[0,5] = 43.39269820653492
$row20[0,6] = 5.275959306258734
$row20[0,7] = 2.632150952014016
$row20[0,8] = 11.21670011910263
$row20[0,9] = 85.08804631142091
$row20[0,10] = 0
$row20[0,11] = 0
$row20[0,12] = 0
$row20[0,13] = 0
$row20[0,14] = 0
$row20[0,15] = 22.28098477411309
$ws.Range("B20:Q20").Value = $row20

$row21 = New-Object "object[,]" 1,16
$row21[0,0] = 30.6775968148569
$row21[0,1] = 0
$row21[0,2] = 2.332840755660176
$row21[0,3] = 39.57372567364283
$row21[0,4] = 31.42443454349621
$row21[0,5] = 46.39810325892545
$row21[0,6] = 5.652473993143921
$row21[0,7] = 2.887686022020177
$row21[0,8] = 11.80885338973472
$row21[0,9] = 90.1856852395535
$row21[0,10] = 0
$row21[0,11] = 0
$row21[0,12] = 0
$row21[0,13] = 0
$row21[0,14] = 0
$row21[0,15] = 23.83243621482351
$ws.Range("B21:Q21").Value = $row21

$row22 = New-Object "object[,]" 1,16
$row22[0,0] = 31.77799867576358
$row22[0,1] = 0
$row22[0,2] = 2.376008484769446
$row22[0,3] = 40.93010792665159
$row22[0,4] = 32.61516141897989
$row22[0,5] = 48.15939008767696
$row22[0,6] = 5.872860374433418
$row22[0,7] = 3.04608517104812
$row22[0,8] = 12.15777632265156
$row22[0,9] = 93.18425791350997
$row22[0,10] = 0
$row22[0,11] = 0
$row22[0,12] = 0
$row22[0,13] = 0
$row22[0,14] = 0
$row22[0,15] = 24.74192583961951
$ws.Range("B22:Q22").Value = $row22

$row23 = New-Object "object[,]" 1,16
$row23[0,0] = 31.19327424365959
$row23[0,1] = 0
$row23[0,2] = 2.35314494378762
$row23[0,3] = 40.21082312578637
$row23[0,4] = 31.98366827679436
$row23[0,5] = 47.22531560716783
$row23[0,6] = 5.755543742063645
$row23[0,7] = 2.95133010413082
$row23[0,8] = 11.97225215802023
$row23[0,9] = 91.38712135191037
$row23[0,10] = 0
$row23[0,11] = 0
$row23[0,12] = 0
$row23[0,13] = 0
$row23[0,14] = 0
$row23[0,15] = 24.25963842951911
$ws.Range("B23:Q23").Value = $row23

$row24 = New-Object "object[,]" 1,16
$row24[0,0] = 28.88424622251212
$row24[0,1] = 0
$row24[0,2] = 2.262714114234934
$row24[0,3] = 37.38415487702425
$row24[0,4] = 29.50294487530883
$row24[0,5] = 43.55563378382656
$row24[0,6] = 5.304189965998743
$row24[0,7] = 2.611260814615723
$row24[0,8] = 11.25435408371167
$row24[0,9] = 84.6651633800618
$row24[0,10] = 0
$row24[0,11] = 0
$row24[0,12] = 0
$row24[0,13] = 0
$row24[0,14] = 0
$row24[0,15] = 22.36376436357785
$ws.Range("B24:Q24").Value = $row24

$row25 = New-Object "object[,]" 1,16
$row25[0,0] = 26.182388850054
$row25[0,1] = 0
$row25[0,2] = 2.157133699284436
$row25[0,3] = 34.11346871640708
$row25[0,4] = 26.67967337940139
$row25[0,5] = 39.30854239459287
$row25[0,6] = 4.800642268843634
$row25[0,7] = 2.742911760467198
$row25[0,8] = 10.44895657665792
$row25[0,9] = 76.8744542187396
$row25[0,10] = 0
$row25[0,11] = 0
$row25[0,12] = 0
$row25[0,13] = 0
$row25[0,14] = 0
$row25[0,15] = 20.32833830186163
$ws.Range("B25:Q25").Value = $row25

